# Apply the recorded edit to the presentation:
#
#  1. Re-style the three tables on slides 14, 15 and 16 to use table
#     style {56ED91A7-8DAA-49C3-BF32-CF3F3A5DD798} instead of
#     {63020C23-DC63-41FF-8568-142B491CE9D0}.
#
#  2. Re-colour the deck's theme palette from the custom "Red Violet"
#     scheme (the "Integral" theme) to the standard Office palette
#     (the colours the "Office Theme" already used elsewhere in the
#     deck), i.e. dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink all move to
#     the stock Office RGB values.

$p = $ppt.ActivePresentation

# --- helper: "RRGGBB" hex string -> COM RGB() integer (0x00BBGGRR) ----
function Convert-HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1. Re-apply table style on every table found on slides 14-16 -----
$newStyleId = "{56ED91A7-8DAA-49C3-BF32-CF3F3A5DD798}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Re-colour the theme (Red Violet -> Office) ---------------------
# MsoThemeColorSchemeIndex order is: dk1, lt1, dk2, lt2,
# accent1..accent6, hyperlink, followed-hyperlink (indices 1-12).
$officeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToRgbInt $officeHex[$i - 1]
}
